$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = 'Background
id="Par1">Medical care workers experienced unprecedented levels of workload and pressure since the outbreak of coronavirus disease 2019 (COVID-19).

 Little is known about its exact impact on medical care workers and related factors in China.

 This study aims to identify the psychological impact of COVID-19 on medical care workers in China.


Methods
id="Par2">From February 23 to March 5, 2020, a cross-sectional survey was conducted among 863 medical care workers from seven provinces in China using standard questionnaires measuring adverse psychological outcomes including Impact of Event Scale-6 (IES-6), Depression, Anxiety and Stress Scale(DASS)and related psychosocial factors like perceived threat, social support and coping strategies.

 Exploratory Factor analysis was performed to identify the dimensions of perceived threat by study participants.

 Multivariate regression was used to examine the determinants of adverse psychological outcomes.


Results
id="Par3">Posttraumatic stress (PTS) were prevalent in this sample of health care professionals, and 40.2% indicated positive screens for significant posttraumatic stress disorder symptoms.

 The proportion of having mild to extremely severe symptoms of depression, anxiety and stress were 13.6, 13.9 and 8.6%, respectively.

 Perceived threat and passive coping strategies were positively correlated to PTS and DASS scores, while perceived social support and active coping strategies were negatively correlated to DASS scores.

 Nurses were more likely to be anxious than others among medical care workers during the COVID-19 epidemic.


Conclusions
id="Par4">Adverse psychological symptoms were prevalent among medical care workers in China during the COVID-19 epidemic.

 Screening for adverse psychological outcomes and developing corresponding preventive measures would be beneficial in decreasing negative psychological outcomes.


'

$ws.Range("E7").Value = '[Ming-Yu%Si%myu_si@163.com%1,  Xiao-You%Su%suxiaoyou@hotmail.com%1,  Yu%Jiang%jiangyu@pumc.edu.cn%2,  Yu%Jiang%jiangyu@pumc.edu.cn%0,  Wen-Jun%Wang%wwjun1973@163.com%1,  Xiao-Fen%Gu%357539328@qq.com%1,  Li%Ma%mali_lele@sina.com%1,  Jing%Li%lijingwcph@scu.edu.cn%0,  Shao-Kai%Zhang%shaokaizhang@126.com%1,  Ze-Fang%Ren%renzef@mail.sysu.edu.cn%1,  Ran%Ren%renran99@163.com%1,  Yuan-Li%Liu%liuyuanli_pumc@163.com%1,  You-Lin%Qiao%qiaoy@cicams.ac.cn%1]'

$ws.Range("D13").Value = 'id="Par1">To assess the psychological effects of the novel coronavirus disease (COVID-19) on medical staff and the general public.

 During the outbreak of COVID-19, an internet-based questionnaire included The Self-rating Depression Scale (SDS), Perceived Stress Scale (PSS-10), and Impact of Event Scale-Revised (IES-R) was used to assess the impact of the pandemic situation on the mental health of medical staff and general population in Wuhan and its surrounding areas.

 Among the 1493 questionnaires completed, 827 (55.39%) of these were men, and 422 (28.27%) of these were medical personnel.

 The results suggest that the outbreak of COVID-19 has affected individuals significantly, the degree of which is related to age, sex, occupation and mental illness.

 There was a significant difference in PSS-10 and IES-R scores between the medical staff and the general population.

 The medical staff showed higher PSS-10 scores (16.813 ± 4.87) and IES-R scores (22.40 ± 12.12) compared to members of the general population PSS-10 (14.80 ± 5.60) and IES-R scores (17.89 ± 13.08).

 However, there was no statistically significant difference between the SDS scores of medical staff (44.52 ± 12.36) and the general public (43.08 ± 11.42).

 In terms of the need for psychological assistance, 50.97% of interviewees responded that they needed psychological counseling, of which medical staff accounted for 65.87% and non-medical staff accounted for 45.10%.

 During the ongoing COVID-19 outbreak, great attention should be paid to the mental health of the population, especially medical staff, and measures such as psychological intervention should be actively carried out for reducing the psychosocial effects.

'

$ws.Range("E13").Value = '[Biao%Chen%NULL%1,  Qing-xian%Li%NULL%1,  Heng%Zhang%NULL%0,  Jia-yong%Zhu%NULL%1,  Xu%Yang%NULL%1,  Yu-hang%Wu%NULL%1,  Jie%Xiong%NULL%2,  Fu%Li%NULL%1,  Hua%Wang%whzkwlh@126.com%0,  Zhi-tao%Chen%chenzhitao830201@163.com%1]'

$ws.Range("D17").Value = 'id="Par1">Coronavirus disease 2019 (COVID-19) is rapidly spreading worldwide, with a staggering number of cases and deaths.

 However, available data on the psychological impacts of COVID-19 on pregnant women are limited.

 The purposes of this study were to assess the prevalence of psychiatric symptoms among pregnant women, and to compare them with non-pregnant women.

 From February 28 to March 12, 2020, a cross-sectional study of pregnant and non-pregnant women was performed in China.

 The online questionnaire was used to collect information of participants.

 The mental health status was assessed by patient health questionnaire, generalized anxiety disorder scale, insomnia severity index, somatization subscale of the symptom checklist 90, and post-traumatic stress disorder (PTSD) checklist-5. Totally, 859 respondents were enrolled, including 544 pregnant women and 315 non-pregnant women.

 In this study, 5.3%, 6.8%, 2.4%, 2.6%, and 0.9% of pregnant women were identified to have symptoms of depression, anxiety, physical discomfort, insomnia, and PTSD, respectively.

 However, the corresponding prevalence rates among non-pregnant women were 17.5%, 17.5%, 2.5%, 5.4%, 5.7%, respectively.

 After adjusting for other covariates, we observed that pregnancy was associated a reduced risk of symptoms of depression (OR = 0.23; 95% CI: 0.12–0.45), anxiety (OR = 0.26; 95% CI: 0.16–0.42), insomnia (OR = 0.19; 95% CI: 0.06–0.58), and PTSD (OR = 0.15; 95% CI: 0.04–0.53) during the COVID-19 epidemic.

 Our results indicate that during the COVID-19 epidemic in China, pregnant women have an advantage of facing mental problems caused by COVID-19, showing fewer depression, anxiety, insomnia, and PTSD symptoms than non-pregnant women.

'

$ws.Range("E17").Value = '[Yongjie%Zhou%NULL%1,  Hui%Shi%NULL%1,  Zhengkui%Liu%NULL%1,  Songxu%Peng%gwxypsx@163.com%2,  Songxu%Peng%gwxypsx@163.com%0,  Ruoxi%Wang%NULL%2,  Ruoxi%Wang%NULL%0,  Ling%Qi%NULL%1,  Zezhi%Li%NULL%1,  Jiezhi%Yang%NULL%1,  Yali%Ren%NULL%1,  Xiuli%Song%NULL%1,  Lingyun%Zeng%NULL%1,  Wei%Qian%NULL%1,  Xiangyang%Zhang%zhangxy@psych.ac.cn%1]'

$ws.Range("D24").Value = 'Background
id="Par1">COVID-19 poses the greatest challenge for the entire world since the Second World War.

 Governments are forced to define strict measures to avoid the spreading of the virus, which may further impose psychological burden for the majority of the population.

 The aim of this study was to evaluate the psychological distress in Austria during the initial stage of the COVID-19 outbreak.


Methods
id="Par2">From 25 March to 3 April 2020, an anonymous online survey was conducted.

 Target group included all members of the Austrian population older than 16 years.

 The survey addressed the following areas (1) and sociodemographic data, (2) physical and mental health; (3) knowledge and concerns about COVID-19; (4) contact with infected people; (5) prevention efforts; (6) need for further information.

 The Impact of Event Scale-Revised (IES-R) and the Depression, Anxiety and Stress Scale (DASS-21) were used to assess mental health.

 Analyses were based on 4126 individuals (74% female, age: M = 38.68, SD = 13.36).


Results
id="Par3">43.3% rated the psychological impact as moderate (5.6%) or severe (37.7%).

 26.5% reported moderate (13.3%) to severe (13.2%) depression; 20.3% moderate (8.9%) to severe (11.4%) anxiety and 21.2% reported to suffer from moderate (10.5%) or severe stress (10.7%).

 Being female, higher age, lower levels of education, concern about family members, internet as main source of information, student or pupil status, poor self-rated health, and downplaying the seriousness of the problem were significantly associated with higher psychological burden.

 Protective factors were the possibility to work in home office, frequent (indirect) contact with family or friends, the availability of virus-specific information, confidence in the diagnosis capability, and physical activity during the crisis.


Conclusion
id="Par4">This study is among the first in Europe on the psychological correlates of the COVID-19 pandemic.

 37.7% of the Austrian study population reported a severe psychological impact on the event and 1 in 10 is considered to suffer from severe depression, anxiety or stress.

 The present findings inform about the identification of protective factors, psychologically vulnerable groups and may guide the development of psychological interventions.


'

$ws.Range("E24").Value = '[Claudia%Traunmüller%claudia.traunmueller@uni-graz.at%1,  Rene%Stefitz%NULL%1,  Kerstin%Gaisbachgrabner%NULL%1,  Andreas%Schwerdtfeger%NULL%1]'

$ws.Range("D54").Value = 'id="Par26">The psychological impact of outbreaks on individuals includes an intense and wide range of psychiatric morbidities.

 People are likely to experience feelings as; worry about being infected or getting sick, increased self-blame, and helplessness.

 This study aimed to assess the impact of COVID-19 on mental health and social support among Egyptian adults during the period of the pandemic.

 This is a cross-sectional observational study using an anonymous online questionnaire.

 The survey was conducted through a link shared on social networking sites.

 It was conducted from 2 May 2020 to 9 May 2020. The general populations of the Egyptian adults were included by using convenience and snowball sampling technique (510 adults).

 Impact Event scale mean 34.3 ± 15. About 211 (41.4%) suffered a severe impact.

 There was an increase in stress from work in 174 (34.1%), financial stress in 284 (55.7%), and stress from home in 320 (62.7%).

 Half of them felt horrified and helpless in 275 (53.9%), and 265 (52%) respectively, while 338 (66.3%) felt apprehensive.

 only 24.2% reported increased support from friends, while increased support from family members in 207 (40.6%).

 46.5% shared their feelings with family members, while 176 (34.5%) shared with others.

 Caring for family members’ feelings increased in 330 (64.7%).

 Age and rural residency were negative predictors for the impact of event score, while female gender or presence of chronic condition was a positive predictor for the impact of event score.

 Covid-19 pandemic has a great psychological impact on adult Egyptians and affected social support.

'

$ws.Range("E54").Value = '[Safaa M.%El-Zoghby%safaaelzoghby86@hotmail.com%1,  Enayat M.%Soltan%enayatsoltan@gmail.com%2,  Enayat M.%Soltan%enayatsoltan@gmail.com%0,  Hend M.%Salama%hind_mikhail@yahoo.com%2,  Hend M.%Salama%hind_mikhail@yahoo.com%0]'

$ws.Range("D66").Value = 'id="Par8">The purposes of this study was to assess the youth mental health after the coronavirus disease 19 (COVID-19) occurred in China two weeks later, and to investigate factors of mental health among youth groups.

 A cross-sectional study was conducted two weeks after the occurrence of COVID-19 in China.

 A total of 584 youth enrolled in this study and completed the question about cognitive status of COVID-19, the General Health Questionnaire(GHQ-12), the PTSD Checklist-Civilian Version (PCL-C) and the Negative coping styles scale.

 Univariate analysis and univariate logistic regression were used to evaluate the effect of COVID-19 on youth mental health.

 The results of this cross-sectional study suggest that nearly 40.4% the sampled youth were found to be prone to psychological problems and 14.4% the sampled youth with Post-traumatic stress disorder (PTSD) symptoms.

 Univariate logistic regression revealed that youth mental health was significantly related to being less educated (OR = 8.71, 95%CI:1.97–38.43), being the enterprise employee (OR = 2.36, 95%CI:1.09–5.09), suffering from the PTSD symptom (OR = 1.05, 95%CI:1.03–1.07) and using negative coping styles (OR = 1.03, 95%CI:1.00–1.07).

 Results of this study suggest that nearly 40.4% of the youth group had a tendency to have psychological problems.

 Thus, this was a remarkable evidence that infectious diseases, such as COVID-19, may have an immense influence on youth mental health.

 Therefor, local governments should develop effective psychological interventions for youth groups, moreover, it is important to consider the educational level and occupation of the youth during the interventions.

'

$ws.Range("E66").Value = '[Leilei%Liang%NULL%1,  Hui%Ren%NULL%1,  Ruilin%Cao%NULL%1,  Yueyang%Hu%NULL%1,  Zeying%Qin%NULL%1,  Chuanen%Li%NULL%1,  Songli%Mei%meisongli@sina.com%1]'
